$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift the "Date" column (C) forward by 63 days.
#    Rows 2-6  (T0 group):      22-May-2024 (45434) -> 24-Jul-2024 (45497)
#    Rows 7-46 (remaining grp): 25-May-2024 (45437) -> 27-Jul-2024 (45500)
# ---------------------------------------------------------------------------
$ws.Range("C2:C6").Value = 45497
$ws.Range("C7:C46").Value = 45500

# ---------------------------------------------------------------------------
# 2. Re-format the Date column from "d-mmm-yy" to a short date (m/d/yyyy),
#    and give the new column L (added below) the same date format.
# ---------------------------------------------------------------------------
$ws.Range("C2:C46").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# 3. Header cell C1 loses its border/shaded style (back to the default,
#    un-styled look) now that the column itself no longer carries a style.
# ---------------------------------------------------------------------------
$ws.Range("C1").Style = "Normal"

# ---------------------------------------------------------------------------
# 4. Add a new (currently empty) column L, pre-formatted the same way as
#    the Date column, for every data row (2-46).
# ---------------------------------------------------------------------------
$ws.Range("L2:L46").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# 5. Recompute the displayed column widths for C and L now that their
#    content/format changed (best-fit, like Excel's auto-fit-column-width).
# ---------------------------------------------------------------------------
$ws.Columns("C").AutoFit() | Out-Null
$ws.Columns("L").AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 6. Match the final on-screen selection (columns J:L were selected).
# ---------------------------------------------------------------------------
$ws.Range("J1:L1048576").Select() | Out-Null
